# Weekly update of the "Ají" (chili pepper) price records reported by the
# Vega Monumental Concepción market. Two new weekly price observations are
# inserted into the historical series (new rows 13 and 19), which pushes the
# previously existing rows down and extends the used range to A1:R37.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 13 for the latest observation, shifting the
# historical rows (old rows 13-35) down by one.
$ws.Rows.Item(13).Insert()
# Row 13
$ws.Cells.Item(13, 1).Value = 11
$ws.Cells.Item(13, 2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(13, 3).Value = 'Bíobío'
$ws.Cells.Item(13, 4).Value = 44425
$ws.Cells.Item(13, 5).Value = 8
$ws.Cells.Item(13, 6).Value = 100112021
$ws.Cells.Item(13, 7).Value = 'Ají'
$ws.Cells.Item(13, 8).Value = 'Inferno'
$ws.Cells.Item(13, 9).Value = 'Primera'
$ws.Cells.Item(13, 10).Value = 30
$ws.Cells.Item(13, 11).Value = 36000
$ws.Cells.Item(13, 12).Value = 38000
$ws.Cells.Item(13, 13).Value = 36667
$ws.Cells.Item(13, 14).Value = '$/caja 25 kilos'
$ws.Cells.Item(13, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(13, 16).Value = 1467
$ws.Cells.Item(13, 17).Value = 25
$ws.Cells.Item(13, 18).Value = 'Hortaliza'

# Insert a second new row at position 19 (after the old row 17, now at row 18)
# for another new observation, shifting the remaining historical rows down once more.
$ws.Rows.Item(19).Insert()
# Row 19
$ws.Cells.Item(19, 1).Value = 11
$ws.Cells.Item(19, 2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(19, 3).Value = 'Bíobío'
$ws.Cells.Item(19, 4).Value = 44427
$ws.Cells.Item(19, 5).Value = 8
$ws.Cells.Item(19, 6).Value = 100112021
$ws.Cells.Item(19, 7).Value = 'Ají'
$ws.Cells.Item(19, 8).Value = 'Americana (o)'
$ws.Cells.Item(19, 9).Value = 'Primera'
$ws.Cells.Item(19, 10).Value = 40
$ws.Cells.Item(19, 11).Value = 75000
$ws.Cells.Item(19, 12).Value = 76000
$ws.Cells.Item(19, 13).Value = 75500
$ws.Cells.Item(19, 14).Value = '$/caja 25 kilos'
$ws.Cells.Item(19, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(19, 16).Value = 3020
$ws.Cells.Item(19, 17).Value = 25
$ws.Cells.Item(19, 18).Value = 'Hortaliza'

